# pedidos.xlsx update — bot appended a new order and authorized a pending one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 12 (pedido 250219_0010): moved from ORÇAMENTO to AUTORIZADO, and
#     the bot stamped the moment the order was placed (data_pedido / col V).
$ws.Range("T12").Value = "AUTORIZADO"
$ws.Range("V12").Value = "2025-02-19 22:34:44"

# --- New row 15: a fresh order (250220_0001) that was CANCELADO.
$ws.Range("A15").Value = "250220_0001"
$ws.Range("B15").Value = "250220_0001_001"
$ws.Range("C15").Value = 6
$ws.Range("D15").Value = "CLIENTE 6"
$ws.Range("E15").Value = "VERDE"
$ws.Range("F15").Value = 29
$ws.Range("G15").Value = "BOX PADRÃO - FIXO - 750MM"
$ws.Range("H15").Value = 20
$ws.Range("I15").Value = "VERDE DE 08MM TEMPERADO"
$ws.Range("J15").Value = "Peça Principal"
$ws.Range("K15").Value = 1
$ws.Range("L15").Value = 1845
$ws.Range("M15").Value = 750
$ws.Range("N15").Value = 1845
$ws.Range("O15").Value = 750
$ws.Range("P15").Value = 1.5
$ws.Range("Q15").Value = 332.75
$ws.Range("R15").Value = 499.12

# nome_pedido (col S) is always stored as text in this sheet, even when the
# value looks numeric (see "48699", "6548", "456" above) - force text entry,
# then drop the leftover number-format so the cell keeps the default style.
$ws.Range("S15").NumberFormat = "@"
$ws.Range("S15").Value = "465"
$ws.Range("S15").ClearFormats()

$ws.Range("T15").Value = "CANCELADO"
$ws.Range("U15").Value = "2025-02-20 00:21:19"

# --- View state: the saved workbook had scrolled right and selected O9.
$ws.Activate()
$ws.Range("O9").Select()
$excel.ActiveWindow.ScrollColumn = 7
$excel.ActiveWindow.ScrollRow = 1
